$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2513.7273
$ws.Range("I28").Value = 697.5
$ws.Range("J28").Value = 7357
$ws.Range("K28").Value = 697.5
$ws.Range("L28").Value = 7357
$ws.Range("M28").Value = -212.5
$ws.Range("N28").Value = -8327
$ws.Range("H32").Value = 986.75
$ws.Range("J32").Value = 1049
$ws.Range("L32").Value = 1049
$ws.Range("N32").Value = -1701
$ws.Range("H40").Value = 5968.8
$ws.Range("I40").Value = 6461
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 6461
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -6286
$ws.Range("N40").Value = -4350
$ws.Range("H53").Value = 167.83333
$ws.Range("I53").Value = 127
$ws.Range("K53").Value = 127
$ws.Range("M53").Value = 510
$ws.Range("H107").Value = 478.7
$ws.Range("I107").Value = 339.29413
$ws.Range("K107").Value = 339.29413
$ws.Range("M107").Value = 1580.70587
$ws.Range("H112").Value = 2099.0454
$ws.Range("J112").Value = 2153.95
$ws.Range("L112").Value = 6461.849999999999
$ws.Range("N112").Value = -8677.849999999999
$ws.Range("H125").Value = 1465.6666
$ws.Range("I125").Value = 1449
$ws.Range("K125").Value = 13041
$ws.Range("M125").Value = -10581
$ws.Range("H141").Value = 2932.5625
$ws.Range("I141").Value = 2910.0833
$ws.Range("K141").Value = 8730.249899999999
$ws.Range("M141").Value = -3550.249899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 15000
$ws.Range("J23").Value = 15000
$ws.Range("L23").Value = 15000
$ws.Range("N23").Value = -15518
$ws.Range("H63").Value = 5370.1333
$ws.Range("I63").Value = 4337.75
$ws.Range("K63").Value = 4337.75
$ws.Range("M63").Value = -3651.75
$ws.Range("H66").Value = 5370.1333
$ws.Range("I66").Value = 4337.75
$ws.Range("K66").Value = 21688.75
$ws.Range("M66").Value = -18256.75
$ws.Range("H97").Value = 854.32355
$ws.Range("I97").Value = 708.4815
$ws.Range("J97").Value = 1416.8572
$ws.Range("K97").Value = 708.4815
$ws.Range("L97").Value = 1416.8572
$ws.Range("M97").Value = -212.4815
$ws.Range("N97").Value = -2408.8572
$ws.Range("H102").Value = 6005.727
$ws.Range("I102").Value = 5896
$ws.Range("J102").Value = 6499.5
$ws.Range("K102").Value = 5896
$ws.Range("L102").Value = 6499.5
$ws.Range("M102").Value = -4274
$ws.Range("N102").Value = -9743.5
$ws.Range("H110").Value = 5094.375
$ws.Range("I110").Value = 1890
$ws.Range("J110").Value = 6162.5
$ws.Range("K110").Value = 1890
$ws.Range("L110").Value = 6162.5
$ws.Range("M110").Value = 155
$ws.Range("N110").Value = -10252.5
$ws.Range("H122").Value = 3460.1052
$ws.Range("I122").Value = 1542.3077
$ws.Range("K122").Value = 4626.9231
$ws.Range("M122").Value = -2176.9231
$ws.Range("H132").Value = 3025.2917
$ws.Range("I132").Value = 1840.4667
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 5521.4001
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -2991.4001
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 22732232
$ws.Range("I20").Value = 35720224
$ws.Range("J20").Value = 3242.625
$ws.Range("K20").Value = 35720224
$ws.Range("L20").Value = 3242.625
$ws.Range("M20").Value = -35719977
$ws.Range("N20").Value = -3736.625
$ws.Range("H86").Value = 1865.08
$ws.Range("I86").Value = 1584.8235
$ws.Range("K86").Value = 1584.8235
$ws.Range("M86").Value = -461.8235
$ws.Range("H89").Value = 1865.08
$ws.Range("I89").Value = 1584.8235
$ws.Range("K89").Value = 7924.1175
$ws.Range("M89").Value = -2308.1175
$ws.Range("H94").Value = 71429576
$ws.Range("I94").Value = 95238440
$ws.Range("J94").Value = 2976
$ws.Range("K94").Value = 95238440
$ws.Range("L94").Value = 2976
$ws.Range("M94").Value = -95237989
$ws.Range("N94").Value = -3878
$ws.Range("H134").Value = 3423.4443
$ws.Range("I134").Value = 3175.147
$ws.Range("K134").Value = 9525.440999999999
$ws.Range("M134").Value = -6990.440999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3540.7354
$ws.Range("I31").Value = 2149.5
$ws.Range("K31").Value = 2149.5
$ws.Range("M31").Value = -1854.5
$ws.Range("H34").Value = 3540.7354
$ws.Range("I34").Value = 2149.5
$ws.Range("K34").Value = 2149.5
$ws.Range("M34").Value = -1947.5
$ws.Range("H92").Value = 31162
$ws.Range("J92").Value = 31162
$ws.Range("L92").Value = 31162
$ws.Range("N92").Value = -36154
$ws.Range("H94").Value = 1875.7646
$ws.Range("I94").Value = 2013
$ws.Range("J94").Value = 1818.5834
$ws.Range("K94").Value = 2013
$ws.Range("L94").Value = 1818.5834
$ws.Range("M94").Value = -1562
$ws.Range("N94").Value = -2720.5834
$ws.Range("H107").Value = 746.2143
$ws.Range("I107").Value = 314.2857
$ws.Range("J107").Value = 1178.1428
$ws.Range("K107").Value = 314.2857
$ws.Range("L107").Value = 1178.1428
$ws.Range("M107").Value = 1605.7143
$ws.Range("N107").Value = -5018.1428
$ws.Range("H122").Value = 4440.029
$ws.Range("I122").Value = 3603.2
$ws.Range("J122").Value = 5067.65
$ws.Range("K122").Value = 10809.6
$ws.Range("L122").Value = 15202.95
$ws.Range("M122").Value = -8359.599999999999
$ws.Range("N122").Value = -20102.95
$ws.Range("H134").Value = 1940.625
$ws.Range("I134").Value = 1724.4186
$ws.Range("J134").Value = 3800
$ws.Range("K134").Value = 5173.2558
$ws.Range("L134").Value = 11400
$ws.Range("M134").Value = -2638.2558
$ws.Range("N134").Value = -16470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 764.1667
$ws.Range("I8").Value = 764.1667
$ws.Range("K8").Value = 2292.5001
$ws.Range("M8").Value = -2153.5001
$ws.Range("H16").Value = 829.8
$ws.Range("I16").Value = 333
$ws.Range("J16").Value = 1575
$ws.Range("K16").Value = 999
$ws.Range("L16").Value = 4725
$ws.Range("M16").Value = -826
$ws.Range("N16").Value = -5071
$ws.Range("H75").Value = 1971.6
$ws.Range("J75").Value = 1971.6
$ws.Range("L75").Value = 5914.799999999999
$ws.Range("N75").Value = -7910.799999999999
$ws.Range("H78").Value = 1971.6
$ws.Range("J78").Value = 1971.6
$ws.Range("L78").Value = 17744.4
$ws.Range("N78").Value = -27728.4
$ws.Range("H107").Value = 577.3
$ws.Range("J107").Value = 446.625
$ws.Range("L107").Value = 1339.875
$ws.Range("N107").Value = -5179.875
$ws.Range("H140").Value = 15184.814
$ws.Range("I140").Value = 11764.117
$ws.Range("K140").Value = 35292.351
$ws.Range("M140").Value = -30112.351

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2848.8572
$ws.Range("I97").Value = 1657
$ws.Range("J97").Value = 10000
$ws.Range("K97").Value = 1657
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = -1161
$ws.Range("N97").Value = -10992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 748.7368
$ws.Range("I16").Value = 754.4706
$ws.Range("K16").Value = 754.4706
$ws.Range("M16").Value = -584.4706
$ws.Range("H61").Value = 7880.8823
$ws.Range("I61").Value = 1701.75
$ws.Range("J61").Value = 22710.8
$ws.Range("K61").Value = 1701.75
$ws.Range("L61").Value = 22710.8
$ws.Range("M61").Value = -1499.75
$ws.Range("N61").Value = -23114.8
$ws.Range("H68").Value = 3144.4
$ws.Range("I68").Value = 3241.1667
$ws.Range("K68").Value = 3241.1667
$ws.Range("M68").Value = -2492.1667
$ws.Range("H71").Value = 3144.4
$ws.Range("I71").Value = 3241.1667
$ws.Range("K71").Value = 16205.8335
$ws.Range("M71").Value = -12461.8335
$ws.Range("H82").Value = 2564.3
$ws.Range("I82").Value = 2549.125
$ws.Range("J82").Value = 2625
$ws.Range("K82").Value = 2549.125
$ws.Range("L82").Value = 2625
$ws.Range("M82").Value = -2188.125
$ws.Range("N82").Value = -3347
$ws.Range("H85").Value = 2564.3
$ws.Range("I85").Value = 2549.125
$ws.Range("J85").Value = 2625
$ws.Range("K85").Value = 2549.125
$ws.Range("L85").Value = 2625
$ws.Range("M85").Value = -1301.125
$ws.Range("N85").Value = -5121
$ws.Range("H113").Value = 7880.8823
$ws.Range("I113").Value = 1701.75
$ws.Range("J113").Value = 22710.8
$ws.Range("K113").Value = 1701.75
$ws.Range("L113").Value = 22710.8
$ws.Range("M113").Value = 468.25
$ws.Range("N113").Value = -27050.8
$ws.Range("H133").Value = 72399.44500000001
$ws.Range("J133").Value = 72399.44500000001
$ws.Range("L133").Value = 72399.44500000001
$ws.Range("N133").Value = -77459.44500000001
$ws.Range("H136").Value = 3457.6191
$ws.Range("I136").Value = 3459.4707
$ws.Range("J136").Value = 3449.75
$ws.Range("K136").Value = 10378.4121
$ws.Range("L136").Value = 10349.25
$ws.Range("M136").Value = -7828.4121
$ws.Range("N136").Value = -15449.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 438.85715
$ws.Range("I113").Value = 254.33333
$ws.Range("K113").Value = 762.99999
$ws.Range("M113").Value = 1407.00001
$ws.Range("H132").Value = 4935.9375
$ws.Range("I132").Value = 4452.273
$ws.Range("K132").Value = 13356.819
$ws.Range("M132").Value = -10826.819
